$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-09-02 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-03 Sunday", 2) | Out-Null

$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "15+34="
$t.Cell(1, 2).Range.Text = "99-71="
$t.Cell(1, 3).Range.Text = "13+55="
$t.Cell(1, 4).Range.Text = "10+78="
$t.Cell(1, 5).Range.Text = "31+45="
$t.Cell(2, 1).Range.Text = "72-52="
$t.Cell(2, 2).Range.Text = "96-90="
$t.Cell(2, 3).Range.Text = "0+78="
$t.Cell(2, 4).Range.Text = "70-47="
$t.Cell(2, 5).Range.Text = "29+61="
$t.Cell(3, 1).Range.Text = "57-44="
$t.Cell(3, 2).Range.Text = "71+26="
$t.Cell(3, 3).Range.Text = "69-25="
$t.Cell(3, 4).Range.Text = "68-31="
$t.Cell(3, 5).Range.Text = "83-38="
$t.Cell(4, 1).Range.Text = "71-21="
$t.Cell(4, 2).Range.Text = "92-5="
$t.Cell(4, 3).Range.Text = "74-46="
$t.Cell(4, 4).Range.Text = "68-20="
$t.Cell(4, 5).Range.Text = "55-19="
$t.Cell(5, 1).Range.Text = "25+3="
$t.Cell(5, 2).Range.Text = "92-10="
$t.Cell(5, 3).Range.Text = "19+21="
$t.Cell(5, 4).Range.Text = "81+3="
$t.Cell(5, 5).Range.Text = "98-18="
$t.Cell(6, 1).Range.Text = "23+60="
$t.Cell(6, 2).Range.Text = "32+21="
$t.Cell(6, 3).Range.Text = "75-68="
$t.Cell(6, 4).Range.Text = "6+50="
$t.Cell(6, 5).Range.Text = "31-2="
$t.Cell(7, 1).Range.Text = "83+1="
$t.Cell(7, 2).Range.Text = "66+4="
$t.Cell(7, 3).Range.Text = "41-6="
$t.Cell(7, 4).Range.Text = "55-3="
$t.Cell(7, 5).Range.Text = "57+42="
$t.Cell(8, 1).Range.Text = "57-22="
$t.Cell(8, 2).Range.Text = "85+14="
$t.Cell(8, 3).Range.Text = "65-27="
$t.Cell(8, 4).Range.Text = "60+24="
$t.Cell(8, 5).Range.Text = "19+49="
$t.Cell(9, 1).Range.Text = "20+21="
$t.Cell(9, 2).Range.Text = "93-43="
$t.Cell(9, 3).Range.Text = "87-70="
$t.Cell(9, 4).Range.Text = "4+91="
$t.Cell(9, 5).Range.Text = "72+6="
$t.Cell(10, 1).Range.Text = "91-6="
$t.Cell(10, 2).Range.Text = "30+27="
$t.Cell(10, 3).Range.Text = "59+8="
$t.Cell(10, 4).Range.Text = "10+66="
$t.Cell(10, 5).Range.Text = "73-70="
$t.Cell(11, 1).Range.Text = "45-0="
$t.Cell(11, 2).Range.Text = "90-49="
$t.Cell(11, 3).Range.Text = "99-15="
$t.Cell(11, 4).Range.Text = "1+35="
$t.Cell(11, 5).Range.Text = "59+18="
$t.Cell(12, 1).Range.Text = "77-75="
$t.Cell(12, 2).Range.Text = "89-71="
$t.Cell(12, 3).Range.Text = "28+6="
$t.Cell(12, 4).Range.Text = "84-41="
$t.Cell(12, 5).Range.Text = "32-32="
$t.Cell(13, 1).Range.Text = "72-23="
$t.Cell(13, 2).Range.Text = "43+32="
$t.Cell(13, 3).Range.Text = "36+38="
$t.Cell(13, 4).Range.Text = "1+47="
$t.Cell(13, 5).Range.Text = "77+15="
$t.Cell(14, 1).Range.Text = "97-13="
$t.Cell(14, 2).Range.Text = "34-14="
$t.Cell(14, 3).Range.Text = "29-1="
$t.Cell(14, 4).Range.Text = "47-43="
$t.Cell(14, 5).Range.Text = "44+15="
$t.Cell(15, 1).Range.Text = "34+19="
$t.Cell(15, 2).Range.Text = "41+13="
$t.Cell(15, 3).Range.Text = "75+15="
$t.Cell(15, 4).Range.Text = "7+53="
$t.Cell(15, 5).Range.Text = "77+15="
$t.Cell(16, 1).Range.Text = "1+48="
$t.Cell(16, 2).Range.Text = "68+7="
$t.Cell(16, 3).Range.Text = "20+65="
$t.Cell(16, 4).Range.Text = "55+21="
$t.Cell(16, 5).Range.Text = "1+95="
$t.Cell(17, 1).Range.Text = "30+58="
$t.Cell(17, 2).Range.Text = "18-4="
$t.Cell(17, 3).Range.Text = "91-41="
$t.Cell(17, 4).Range.Text = "61+11="
$t.Cell(17, 5).Range.Text = "95-67="
$t.Cell(18, 1).Range.Text = "19-18="
$t.Cell(18, 2).Range.Text = "90-40="
$t.Cell(18, 3).Range.Text = "48+14="
$t.Cell(18, 4).Range.Text = "15+73="
$t.Cell(18, 5).Range.Text = "73+14="
$t.Cell(19, 1).Range.Text = "32+56="
$t.Cell(19, 2).Range.Text = "49+46="
$t.Cell(19, 3).Range.Text = "11+84="
$t.Cell(19, 4).Range.Text = "43-37="
$t.Cell(19, 5).Range.Text = "39+25="
$t.Cell(20, 1).Range.Text = "32-5="
$t.Cell(20, 2).Range.Text = "25+18="
$t.Cell(20, 3).Range.Text = "22+39="
$t.Cell(20, 4).Range.Text = "71-64="
$t.Cell(20, 5).Range.Text = "30+69="
